$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 61: new "8ème semaine" section header, styled like the other section headers (e.g. row 9)
$ws.Range("A9:C9").Copy()
$ws.Range("A61:C61").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A61").Value2 = "8ème semaine "
$ws.Range("A61:C61").Merge()

# Prepare rows 62 and 63 formats first (so the shared-string pool still
# receives the two new strings in the same order the original author typed
# them in: row 63's text was entered before row 62's).
$ws.Range("A49:C49").Copy()
$ws.Range("A62:C62").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A50:C50").Copy()
$ws.Range("A63:C63").PasteSpecial(-4122) # xlPasteFormats

# Row 63: activity entry - date 2018-04-17 (serial 43207), wrapped text in B, plain text in C
$ws.Range("A63").Value2 = 43207
$ws.Range("B63").Value2 = "J'ai changé le nom de ma base de données, j'ai géré les injections sql pour le site."
$ws.Range("C63").Value2 = "2 périodes"

# Row 62: activity entry - date 2018-04-17 (serial 43207), wrapped text in B, plain text in C
$ws.Range("A62").Value2 = 43207
$ws.Range("B62").Value2 = "J'ai rempli ma documentation, j'ai fait ma conclusion et j'ai ajouté des choses dans mon tableau de bord. J'ai ajouté aussi des choses dans le journal de bord."
$ws.Range("C62").Value2 = "3 périodes"
$ws.Rows(62).RowHeight = 30

$excel.CutCopyMode = 0

# --- View state updates ---
$ws.Application.ActiveWindow.Zoom = 125
[void]$ws.Range("C63").Select()
